$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.534.97"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "3.384.02"
$ws.Range("E3").Value = "  -1.23%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'405.77"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").Value = "'126.14"
$ws.Range("E6").Value = "  -3.08%  "

$ws.Range("D7").Value = "'0.609"
$ws.Range("E7").Value = "  -3.13%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.708"
$ws.Range("E9").Value = "  -4.26%  "

$ws.Range("D10").Value = "'0.130"
$ws.Range("E10").Value = "  -12.11%  "

$ws.Range("D11").Value = "'41.57"
$ws.Range("E11").Value = "  -2.74%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.140"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.919.08"
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").Value = "'8.90"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'20.11"
$ws.Range("E15").Value = "  -4.73%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000199"
$ws.Range("E16").Value = "  -13.93%  "

$ws.Range("D17").Value = "3.385.00"
$ws.Range("E17").Value = "  -1.35%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.01"
$ws.Range("E18").Value = "  -3.46%  "

$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D19").Value = "'1.05"
$ws.Range("E19").Value = "  -2.18%  "

$ws.Range("D20").Value = "61.513.79"
$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").Value = "'473.70"
$ws.Range("E21").Value = "  +19.04%  "

$ws.Range("D22").Value = "'88.17"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").Value = "'3.16"
$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("D24").Value = "'12.88"
$ws.Range("E24").Value = "  -2.84%  "

$ws.Range("D25").Value = "'3.23"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").Value = "'32.75"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("E27").Value = "  +3.66%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "'7.84"
$ws.Range("E29").Value = "  +3.26%  "

$ws.Range("D30").Value = "'2.60"
$ws.Range("E30").Value = "  -4.42%  "

$ws.Range("D31").Value = "'11.61"
$ws.Range("E31").Value = "  -2.62%  "

$ws.Range("D32").Value = "'0.164"
$ws.Range("E32").Value = "  -4.61%  "

$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = "  -7.66%  "

$ws.Range("D34").Value = "'40.14"
$ws.Range("E34").Value = "  -6.89%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "'54.81"
$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").Value = "'0.0476"
$ws.Range("E37").Value = "  -5.39%  "

$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("D39").Value = "'147.36"
$ws.Range("E39").Value = "  +4.13%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.322"
$ws.Range("E40").Value = "  +3.18%  "

$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value = "'3.30"
$ws.Range("E41").Value = "  -2.92%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.131"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("D43").Value = "'2.87"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").Value = "'2.03"
$ws.Range("E44").Value = "  +2.36%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.54"
$ws.Range("E45").Value = "  +4.00%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'4.09"
$ws.Range("E46").Value = "  -0.64%  "

$ws.Range("D47").Value = "'2.30"
$ws.Range("E47").Value = "  +16.58%  "

$ws.Range("D48").Value = "'16.06"
$ws.Range("E48").Value = "  -4.23%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'21.61"
$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.142"
$ws.Range("E50").Value = "  +7.92%  "

$ws.Range("D51").Value = "'111.86"
$ws.Range("E51").Value = "  +14.87%  "
